# Adds the "INTRODUCCIÓN" intro paragraph described in the commit:
#   "feat: Agregado parrafo de introducción - Barreto"
#
# Before: the document ends with a heading paragraph ("INTRODUCCIÓN")
#         that also carries the hidden "_GoBack" bookmark right before
#         its run.
# After:  a new body paragraph is appended after the heading with the
#         introductory text, and the hidden "_GoBack" bookmark is moved
#         into that new paragraph (splitting its text into two runs
#         around the bookmark).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The existing "_GoBack" bookmark currently sits right before the
#    "INTRODUCCIÓN" run. It needs to move into the new paragraph, so
#    delete it here and re-add it later at the right spot.
# ------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# ------------------------------------------------------------------
# 2. "INTRODUCCIÓN" is the last paragraph in the document body.
#    Insert a brand-new paragraph right after it.
# ------------------------------------------------------------------
$introPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$introRange = $introPara.Range
$introRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Type the full sentence first (while the paragraph still inherits the
# heading's formatting) so the run does not pick up bold text.
$nr = $newPara.Range
$nr.Text = "El problema identificado es la inexistencia de unas bases de trabajo automatizado que asista a los farmacéuticos para un mejor control del negocio."

# Now switch the paragraph to the "Prrafodelista" (List Paragraph) style
# used by the rest of the body text (drops the inherited numbering too)
# and size it like normal body copy (12pt / sz 24).
$newPara.Style = "Prrafodelista"
$nr2 = $newPara.Range
$nr2.Font.Size = 12

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark inside the new paragraph, right
#    after "El problema iden" (splitting the sentence into the two
#    runs "El problema iden" / "tificado es ...").
# ------------------------------------------------------------------
$bmPos = $nr2.Start + 16
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
